$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update values in column P (and a couple in column N) ---
$ws.Range("P3").Value = 82.5
$ws.Range("P4").Value = 100
$ws.Range("P5").Value = 66
$ws.Range("P6").Value = 105
$ws.Range("P7").Value = 86.5
$ws.Range("P8").Value = 51
$ws.Range("P9").Value = 60.5
$ws.Range("P10").Value = 75.5

$ws.Range("N11").Value = 19
$ws.Range("P11").Value = 95.5

$ws.Range("N13").Value = 8.1
$ws.Range("P13").Value = 64

$ws.Range("P14").Value = 50
$ws.Range("P15").Value = 61
$ws.Range("P17").Value = 48.5
$ws.Range("P19").Value = 55.5
$ws.Range("P20").Value = 49.5
$ws.Range("P21").Value = 63.5
$ws.Range("P22").Value = 72
$ws.Range("P23").Value = 63.5
$ws.Range("P24").Value = 105
$ws.Range("P25").Value = 61.5

# --- Update sheet view: zoom level + active-cell selection ---
$ws.Activate()
$window = $excel.ActiveWindow
$window.Zoom = 130
[void]$ws.Range("H9").Select()
